$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-03-08 Saturday" "2025-03-09 Sunday"
Replace-Text "121×6=726" "117×8=936"
Replace-Text "123×6=738" "853×8=6824"
Replace-Text "857×9=7713" "527×3=1581"
Replace-Text "384×5=1920" "972×6=5832"
Replace-Text "273×5=1365" "135×9=1215"
Replace-Text "877×5=4385" "347×3=1041"
Replace-Text "226×5=1130" "465×8=3720"
Replace-Text "702×4=2808" "510×7=3570"
Replace-Text "252×2=504" "405×8=3240"
Replace-Text "680×3=2040" "916×9=8244"
Replace-Text "473×8=3784" "107×2=214"
Replace-Text "239×9=2151" "977×9=8793"
Replace-Text "129×2=258" "830×6=4980"
Replace-Text "726×9=6534" "584×6=3504"
Replace-Text "787×2=1574" "951×6=5706"
Replace-Text "826×7=5782" "354×9=3186"
Replace-Text "797×6=4782" "935×8=7480"
Replace-Text "274×2=548" "401×7=2807"
Replace-Text "803×5=4015" "667×3=2001"
Replace-Text "953×4=3812" "654×2=1308"
Replace-Text "291×3=873" "473×5=2365"
Replace-Text "803×3=2409" "373×5=1865"
Replace-Text "126×9=1134" "215×2=430"
Replace-Text "176×9=1584" "487×6=2922"
Replace-Text "319×8=2552" "359×8=2872"
